$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.026.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.525.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.525.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.135.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000200"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.517.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.009.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "417.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.598"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.670.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.523.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "174.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0818"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.855"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  -7.34%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.71%  "
